$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log row appended after the existing last row (54) -> row 55.
$newRow = 55
$prevRow = 54

# Copy formatting from the previous row so the new row matches the sheet's
# existing style (center/center alignment, same font/border as the rest of
# the log table) without minting unused style entries.
$ws.Range("A$prevRow`:H$prevRow").Copy()
$ws.Range("A$newRow`:H$newRow").PasteSpecial(-4122)

$ws.Range("A$newRow").Value = "2025-08-25 06:50:56 UTC"
$ws.Range("B$newRow").Value = "2025-08-25 12:20:56 IST"
$ws.Range("C$newRow").Value = "SKIPPED"
$ws.Range("D$newRow").Value = "No change in PDF. Skipping download & Excel update."
$ws.Range("E$newRow").Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"
$ws.Range("F$newRow").Value = ""
$ws.Range("G$newRow").Value = 0
$ws.Range("H$newRow").Value = ""

Write-Output "done"
